$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Helper: copy the cell format (style) from a stable source cell onto a
# target cell using Copy + PasteSpecial(xlPasteFormats), so the resulting
# cellXfs index matches the one the real workbook ends up with.
# ---------------------------------------------------------------------------
function Copy-Format($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial(-4122) | Out-Null
}

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# New header cells I1:K1 -> MP / AttackRange / AttackAngle
# I1's style changes (was s=4, now s=2) -- copy from C1 which already has s=2.
# J1 / K1 keep their existing style (s=5), so a plain value assignment is safe.
# ---------------------------------------------------------------------------
Copy-Format 1 3 1 9   # C1 (s=2) -> I1
$ws.Cells.Item(1, 9).Value = "MP"
$ws.Cells.Item(1, 10).Value = "AttackRange"
$ws.Cells.Item(1, 11).Value = "AttackAngle"

# ---------------------------------------------------------------------------
# Row 2 (SkillID 101 - JumpSkill): DelayTime 6 -> 5; new MP/Range/Angle values.
# I2/J2/K2 already carry the correct style (s=4 / s=5 / s=5), so plain value
# assignment is enough.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 45

# ---------------------------------------------------------------------------
# Row 3 (SkillID 102 - StoneSlash): DelayTime 8 -> 6; new MP/Range/Angle.
# Row 3's default style is s=4, which matches I3's target style, and J3/K3
# already carry their correct explicit styles (s=3 / s=6).
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 6).Value = 6
$ws.Cells.Item(3, 9).Value = 3
$ws.Cells.Item(3, 10).Value = 2.4
$ws.Cells.Item(3, 11).Value = 100

# ---------------------------------------------------------------------------
# Row 4 (SkillID 103 - AssassinAttack): DelayTime 20 -> 8; new MP/Range/Angle.
# Same reasoning as row 3: row default s=4 covers I4; J4/K4 already s=3.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(4, 9).Value = 4
$ws.Cells.Item(4, 10).Value = 5
$ws.Cells.Item(4, 11).Value = 90

# ---------------------------------------------------------------------------
# Row 5 (SkillID 104 - StingAttack): DelayTime 25 -> 10; Damage 20 -> 30;
# new MP/Range/Angle. J5/K5 need an explicit style of s=7 (row default is
# s=4), so copy formats from C2 which already carries s=7.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 6).Value = 10
$ws.Cells.Item(5, 7).Value = 30
$ws.Cells.Item(5, 9).Value = 5
Copy-Format 2 3 5 10   # C2 (s=7) -> J5
Copy-Format 2 3 5 11   # C2 (s=7) -> K5
$ws.Cells.Item(5, 10).Value = 2.5
$ws.Cells.Item(5, 11).Value = 90

# ---------------------------------------------------------------------------
# Row 6 (SkillID 201 - Sword2Default): new MP/Range/Angle.
# J6/K6 need explicit style s=5 (row default is s=4) -> copy from J1 (s=5).
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 9).Value = 3
Copy-Format 1 10 6 10   # J1 (s=5) -> J6
Copy-Format 1 10 6 11   # J1 (s=5) -> K6
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(6, 11).Value = 2

# ---------------------------------------------------------------------------
# Row 7 (SkillID 202 - Sword2Skill1): new MP/Range/Angle.
# J7/K7 need explicit style s=3 (row default is s=4) -> copy from J4 (s=3).
# ---------------------------------------------------------------------------
$ws.Cells.Item(7, 9).Value = 3
Copy-Format 4 10 7 10   # J4 (s=3) -> J7
Copy-Format 4 10 7 11   # J4 (s=3) -> K7
$ws.Cells.Item(7, 10).Value = 2.4
$ws.Cells.Item(7, 11).Value = 2.4

# ---------------------------------------------------------------------------
# Row 8 (SkillID 203 - Sword2Skill2): new MP/Range/Angle.
# J8/K8 need explicit style s=3 -> copy from J4 (s=3).
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 9).Value = 3
Copy-Format 4 10 8 10   # J4 (s=3) -> J8
Copy-Format 4 10 8 11   # J4 (s=3) -> K8
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(8, 11).Value = 5

# ---------------------------------------------------------------------------
# Row 9 (SkillID 301 - Sword3Default): new MP/Range/Angle.
# J9/K9 need explicit style s=7 -> copy from C2 (s=7).
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 9).Value = 3
Copy-Format 2 3 9 10   # C2 (s=7) -> J9
Copy-Format 2 3 9 11   # C2 (s=7) -> K9
$ws.Cells.Item(9, 10).Value = 2.5
$ws.Cells.Item(9, 11).Value = 2.5

# ---------------------------------------------------------------------------
# Row 10 (SkillID 302 - Sword3Skill1): new MP/Range/Angle.
# J10/K10 need explicit style s=7 -> copy from C2 (s=7).
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 9).Value = 3
Copy-Format 2 3 10 10   # C2 (s=7) -> J10
Copy-Format 2 3 10 11   # C2 (s=7) -> K10
$ws.Cells.Item(10, 10).Value = 10
$ws.Cells.Item(10, 11).Value = 10

# ---------------------------------------------------------------------------
# Row 11 (SkillID 303 - Sword3Skill2): new MP/Range/Angle.
# J11/K11 need explicit style s=7 -> copy from C2 (s=7).
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 9).Value = 3
Copy-Format 2 3 11 10   # C2 (s=7) -> J11
Copy-Format 2 3 11 11   # C2 (s=7) -> K11
$ws.Cells.Item(11, 10).Value = 10
$ws.Cells.Item(11, 11).Value = 10

# ---------------------------------------------------------------------------
# Selection moves from G6 to H7.
# ---------------------------------------------------------------------------
$ws.Range("H7").Select()
